$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22: cohort 2020, period 5 -> num_customers 33 -> 34, retention_rate recalculated (cohort_size 2654 unchanged)
$ws.Range("C22").Value = 34
$ws.Range("E22").Value = 34 / 2654

# Row 31: cohort 2021, period 3 -> num_customers 61 -> 64, retention_rate recalculated (cohort_size 2312 unchanged)
$ws.Range("C31").Value = 64
$ws.Range("E31").Value = 64 / 2312

# Row 34: cohort 2022, period 2 -> num_customers 93 -> 95, retention_rate recalculated (cohort_size 2256 unchanged)
$ws.Range("C34").Value = 95
$ws.Range("E34").Value = 95 / 2256

# Row 36: cohort 2023, period 1 -> num_customers 153 -> 155, retention_rate recalculated (cohort_size 1930 unchanged)
$ws.Range("C36").Value = 155
$ws.Range("E36").Value = 155 / 1930

# Row 37: cohort 2025, period 0 -> num_customers 1028 -> 1032, cohort_size 1028 -> 1032 (retention_rate stays 1)
$ws.Range("C37").Value = 1032
$ws.Range("D37").Value = 1032
